$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.724.97'
$ws.Range('E2').Value = '  -0.26%  '
$ws.Range('D3').Value = '1.632.98'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.95'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('E6').Value = '  -0.66%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -1.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.63'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.99%  '
$ws.Range('E11').Value = '  +1.03%  '
$ws.Range('E12').Value = '  -0.33%  '
$ws.Range('D13').Value = '1.858.25'
$ws.Range('E13').Value = '  -0.18%  '
$ws.Range('D14').Value = '1.632.49'
$ws.Range('E14').Value = '  -0.06%  '
$ws.Range('E15').Value = '  -0.34%  '
$ws.Range('D16').Value = '0.0₃0763'
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.65'
$ws.Range('D17').ClearFormats()
$ws.Range('D18').Value = '25.746.47'
$ws.Range('E18').Value = '  -0.23%  '
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('E20').Value = '  +1.63%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '193.73'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.94'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('E23').Value = '  +2.36%  '
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('E25').Value = '  +3.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '141.98'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.20%  '
$ws.Range('E27').Value = '  -0.96%  '
$ws.Range('E28').Value = '  +0.97%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.51'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('E30').Value = '  -0.15%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0492'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.44%  '
$ws.Range('E32').Value = '  +0.68%  '
$ws.Range('E33').Value = '  -0.41%  '
$ws.Range('E34').Value = '  +0.41%  '
$ws.Range('E35').Value = '  +0.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.899'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.29%  '
$ws.Range('D37').Value = '1.124.21'
$ws.Range('E37').Value = '  -0.57%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.547'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.58%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.52'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.69%  '
$ws.Range('E40').Value = '  -1.08%  '
$ws.Range('E42').Value = '  +2.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.66'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.78%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.804'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.62%  '
$ws.Range('D45').Value = '1.767.70'
$ws.Range('E45').Value = '  -0.29%  '
$ws.Range('D46').Value = '0.0₆0110'
$ws.Range('E46').Value = '  -3.83%  '
$ws.Range('E47').Value = '  -1.03%  '
$ws.Range('E48').Value = '  -2.07%  '
$ws.Range('B50').Value = 'SynthetixNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.34'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.22%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.53'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -3.11%  '
